# Auto-generated Excel COM-interop edit script
# Applies updated market-price figures across ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 6450   # H51: 6633.3335 -> 6450
$ws.Cells.Item(51, 10).Value = 6450   # J51: 6633.3335 -> 6450
$ws.Cells.Item(51, 12).Value = 6450   # L51: 6633.3335 -> 6450
$ws.Cells.Item(51, 14).Value = -7418   # N51: -7601.3335 -> -7418

$ws.Cells.Item(70, 8).Value = 18999.5   # H70: 18874.125 -> 18999.5
$ws.Cells.Item(70, 9).Value = 21499   # I70: 17998 -> 21499
$ws.Cells.Item(70, 10).Value = 18166.334   # J70: 18999.285 -> 18166.334
$ws.Cells.Item(70, 11).Value = 64497   # K70: 53994 -> 64497
$ws.Cells.Item(70, 12).Value = 54499.00199999999   # L70: 56997.855 -> 54499.00199999999
$ws.Cells.Item(70, 13).Value = -64227   # M70: -53724 -> -64227
$ws.Cells.Item(70, 14).Value = -55039.00199999999   # N70: -57537.855 -> -55039.00199999999

$ws.Cells.Item(73, 8).Value = 18999.5   # H73: 18874.125 -> 18999.5
$ws.Cells.Item(73, 9).Value = 21499   # I73: 17998 -> 21499
$ws.Cells.Item(73, 10).Value = 18166.334   # J73: 18999.285 -> 18166.334
$ws.Cells.Item(73, 11).Value = 64497   # K73: 53994 -> 64497
$ws.Cells.Item(73, 12).Value = 54499.00199999999   # L73: 56997.855 -> 54499.00199999999
$ws.Cells.Item(73, 13).Value = -63561   # M73: -53058 -> -63561
$ws.Cells.Item(73, 14).Value = -56371.00199999999   # N73: -58869.855 -> -56371.00199999999

$ws.Cells.Item(87, 8).Value = 19869.564   # H87: 19833.334 -> 19869.564
$ws.Cells.Item(87, 10).Value = 19869.564   # J87: 19833.334 -> 19869.564
$ws.Cells.Item(87, 12).Value = 19869.564   # L87: 19833.334 -> 19869.564
$ws.Cells.Item(87, 14).Value = -22365.564   # N87: -22329.334 -> -22365.564

$ws.Cells.Item(90, 8).Value = 19869.564   # H90: 19833.334 -> 19869.564
$ws.Cells.Item(90, 10).Value = 19869.564   # J90: 19833.334 -> 19869.564
$ws.Cells.Item(90, 12).Value = 59608.692   # L90: 59500.00199999999 -> 59608.692
$ws.Cells.Item(90, 14).Value = -72088.692   # N90: -71980.00199999999 -> -72088.692

$ws.Cells.Item(92, 8).Value = 5555719.5   # H92: 5555721 -> 5555719.5
$ws.Cells.Item(92, 9).Value = 170.375   # I92: 172.25 -> 170.375
$ws.Cells.Item(92, 11).Value = 170.375   # K92: 172.25 -> 170.375
$ws.Cells.Item(92, 13).Value = 1077.625   # M92: 1075.75 -> 1077.625

$ws.Cells.Item(100, 8).Value = 4476.0557   # H100: 4821.879 -> 4476.0557
$ws.Cells.Item(100, 9).Value = 3184.238   # I100: 3308.95 -> 3184.238
$ws.Cells.Item(100, 10).Value = 6284.6   # J100: 7149.4614 -> 6284.6
$ws.Cells.Item(100, 11).Value = 3184.238   # K100: 3308.95 -> 3184.238
$ws.Cells.Item(100, 12).Value = 6284.6   # L100: 7149.4614 -> 6284.6
$ws.Cells.Item(100, 13).Value = -2643.238   # M100: -2767.95 -> -2643.238
$ws.Cells.Item(100, 14).Value = -7366.6   # N100: -8231.4614 -> -7366.6

$ws.Cells.Item(106, 8).Value = 3912   # H106: 3711.9333 -> 3912
$ws.Cells.Item(106, 9).Value = 3607.2727   # I106: 3389.9167 -> 3607.2727
$ws.Cells.Item(106, 10).Value = 4750   # J106: 5000 -> 4750
$ws.Cells.Item(106, 11).Value = 3607.2727   # K106: 3389.9167 -> 3607.2727
$ws.Cells.Item(106, 12).Value = 4750   # L106: 5000 -> 4750
$ws.Cells.Item(106, 13).Value = -2976.2727   # M106: -2758.9167 -> -2976.2727
$ws.Cells.Item(106, 14).Value = -6012   # N106: -6262 -> -6012

$ws.Cells.Item(107, 8).Value = 444.07693   # H107: 889.5625 -> 444.07693
$ws.Cells.Item(107, 9).Value = 516.3   # I107: 803.3333 -> 516.3
$ws.Cells.Item(107, 10).Value = 203.33333   # J107: 1148.25 -> 203.33333
$ws.Cells.Item(107, 11).Value = 516.3   # K107: 803.3333 -> 516.3
$ws.Cells.Item(107, 12).Value = 203.33333   # L107: 1148.25 -> 203.33333
$ws.Cells.Item(107, 13).Value = 1403.7   # M107: 1116.6667 -> 1403.7
$ws.Cells.Item(107, 14).Value = -4043.33333   # N107: -4988.25 -> -4043.33333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(34, 8).Value = 265014   # H34: 91672 -> 265014
$ws.Cells.Item(34, 10).Value = 265014   # J34: 91672 -> 265014
$ws.Cells.Item(34, 12).Value = 265014   # L34: 91672 -> 265014
$ws.Cells.Item(34, 14).Value = -265556   # N34: -92214 -> -265556

$ws.Cells.Item(62, 8).Value = 27500   # H62: 0 -> 27500
$ws.Cells.Item(62, 10).Value = 27500   # J62: 0 -> 27500
$ws.Cells.Item(62, 12).Value = 27500   # L62: 0 -> 27500
$ws.Cells.Item(62, 14).Value = -28748   # N62: None -> -28748

$ws.Cells.Item(65, 8).Value = 27500   # H65: 0 -> 27500
$ws.Cells.Item(65, 10).Value = 27500   # J65: 0 -> 27500
$ws.Cells.Item(65, 12).Value = 82500   # L65: 0 -> 82500
$ws.Cells.Item(65, 14).Value = -88740   # N65: None -> -88740

$ws.Cells.Item(132, 8).Value = 2318.5   # H132: 2553 -> 2318.5
$ws.Cells.Item(132, 9).Value = 2318.5   # I132: 2553 -> 2318.5
$ws.Cells.Item(132, 11).Value = 6955.5   # K132: 7659 -> 6955.5
$ws.Cells.Item(132, 13).Value = -4425.5   # M132: -5129 -> -4425.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 47407.59   # H86: 2076.8572 -> 47407.59
$ws.Cells.Item(86, 9).Value = 54446.895   # I86: 2121.2307 -> 54446.895
$ws.Cells.Item(86, 10).Value = 2825.3333   # J86: 1500 -> 2825.3333
$ws.Cells.Item(86, 11).Value = 54446.895   # K86: 2121.2307 -> 54446.895
$ws.Cells.Item(86, 12).Value = 2825.3333   # L86: 1500 -> 2825.3333
$ws.Cells.Item(86, 13).Value = -53323.895   # M86: -998.2307000000001 -> -53323.895
$ws.Cells.Item(86, 14).Value = -5071.3333   # N86: -3746 -> -5071.3333

$ws.Cells.Item(89, 8).Value = 47407.59   # H89: 2076.8572 -> 47407.59
$ws.Cells.Item(89, 9).Value = 54446.895   # I89: 2121.2307 -> 54446.895
$ws.Cells.Item(89, 10).Value = 2825.3333   # J89: 1500 -> 2825.3333
$ws.Cells.Item(89, 11).Value = 272234.475   # K89: 10606.1535 -> 272234.475
$ws.Cells.Item(89, 12).Value = 14126.6665   # L89: 7500 -> 14126.6665
$ws.Cells.Item(89, 13).Value = -266618.475   # M89: -4990.1535 -> -266618.475
$ws.Cells.Item(89, 14).Value = -25358.6665   # N89: -18732 -> -25358.6665

$ws.Cells.Item(94, 8).Value = 2586734.2   # H94: 2648331.2 -> 2586734.2
$ws.Cells.Item(94, 9).Value = 2169.5334   # I94: 2255.862 -> 2169.5334
$ws.Cells.Item(94, 11).Value = 2169.5334   # K94: 2255.862 -> 2169.5334
$ws.Cells.Item(94, 13).Value = -1718.5334   # M94: -1804.862 -> -1718.5334

$ws.Cells.Item(96, 8).Value = 20965.166   # H96: 24310.2 -> 20965.166
$ws.Cells.Item(96, 9).Value = 5358.4   # I96: 5638 -> 5358.4
$ws.Cells.Item(96, 11).Value = 5358.4   # K96: 5638 -> 5358.4
$ws.Cells.Item(96, 13).Value = -2612.4   # M96: -2892 -> -2612.4

$ws.Cells.Item(140, 8).Value = 39769.23   # H140: 60000 -> 39769.23
$ws.Cells.Item(140, 10).Value = 39769.23   # J140: 60000 -> 39769.23
$ws.Cells.Item(140, 12).Value = 39769.23   # L140: 60000 -> 39769.23
$ws.Cells.Item(140, 14).Value = -50129.23   # N140: -70360 -> -50129.23

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2936.4482   # H31: 3071.577 -> 2936.4482
$ws.Cells.Item(31, 9).Value = 2325.5454   # I31: 2414 -> 2325.5454
$ws.Cells.Item(31, 11).Value = 2325.5454   # K31: 2414 -> 2325.5454
$ws.Cells.Item(31, 13).Value = -2030.5454   # M31: -2119 -> -2030.5454

$ws.Cells.Item(34, 8).Value = 2936.4482   # H34: 3071.577 -> 2936.4482
$ws.Cells.Item(34, 9).Value = 2325.5454   # I34: 2414 -> 2325.5454
$ws.Cells.Item(34, 11).Value = 2325.5454   # K34: 2414 -> 2325.5454
$ws.Cells.Item(34, 13).Value = -2123.5454   # M34: -2212 -> -2123.5454

$ws.Cells.Item(58, 8).Value = 3150.1428   # H58: 3393.8572 -> 3150.1428
$ws.Cells.Item(58, 9).Value = 2758   # I58: 3073.4 -> 2758
$ws.Cells.Item(58, 10).Value = 3856   # J58: 4195 -> 3856
$ws.Cells.Item(58, 11).Value = 2758   # K58: 3073.4 -> 2758
$ws.Cells.Item(58, 12).Value = 3856   # L58: 4195 -> 3856
$ws.Cells.Item(58, 13).Value = -2555   # M58: -2870.4 -> -2555
$ws.Cells.Item(58, 14).Value = -4262   # N58: -4601 -> -4262

$ws.Cells.Item(99, 8).Value = 2530.25   # H99: 2749.5 -> 2530.25
$ws.Cells.Item(99, 9).Value = 2573.6667   # I99: 2749.5 -> 2573.6667
$ws.Cells.Item(99, 10).Value = 2400   # J99: 0 -> 2400
$ws.Cells.Item(99, 11).Value = 2573.6667   # K99: 2749.5 -> 2573.6667
$ws.Cells.Item(99, 12).Value = 2400   # L99: 0 -> 2400
$ws.Cells.Item(99, 13).Value = -1075.6667   # M99: -1251.5 -> -1075.6667
$ws.Cells.Item(99, 14).Value = -5396   # N99: None -> -5396

$ws.Cells.Item(126, 8).Value = 2530.25   # H126: 2749.5 -> 2530.25
$ws.Cells.Item(126, 9).Value = 2573.6667   # I126: 2749.5 -> 2573.6667
$ws.Cells.Item(126, 10).Value = 2400   # J126: 0 -> 2400
$ws.Cells.Item(126, 11).Value = 7721.000100000001   # K126: 8248.5 -> 7721.000100000001
$ws.Cells.Item(126, 12).Value = 7200   # L126: 0 -> 7200
$ws.Cells.Item(126, 13).Value = -5251.000100000001   # M126: -5778.5 -> -5251.000100000001
$ws.Cells.Item(126, 14).Value = -12140   # N126: None -> -12140

$ws.Cells.Item(127, 8).Value = 30000   # H127: 0 -> 30000
$ws.Cells.Item(127, 10).Value = 30000   # J127: 0 -> 30000
$ws.Cells.Item(127, 12).Value = 30000   # L127: 0 -> 30000
$ws.Cells.Item(127, 14).Value = -39920   # N127: None -> -39920

$ws.Cells.Item(136, 8).Value = 3150.1428   # H136: 3393.8572 -> 3150.1428
$ws.Cells.Item(136, 9).Value = 2758   # I136: 3073.4 -> 2758
$ws.Cells.Item(136, 10).Value = 3856   # J136: 4195 -> 3856
$ws.Cells.Item(136, 11).Value = 8274   # K136: 9220.200000000001 -> 8274
$ws.Cells.Item(136, 12).Value = 11568   # L136: 12585 -> 11568
$ws.Cells.Item(136, 13).Value = -5724   # M136: -6670.200000000001 -> -5724
$ws.Cells.Item(136, 14).Value = -16668   # N136: -17685 -> -16668

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(40, 8).Value = 11695.875   # H40: 13820.875 -> 11695.875
$ws.Cells.Item(40, 9).Value = 10078.75   # I40: 14328.75 -> 10078.75
$ws.Cells.Item(40, 11).Value = 10078.75   # K40: 14328.75 -> 10078.75
$ws.Cells.Item(40, 13).Value = -9927.75   # M40: -14177.75 -> -9927.75

$ws.Cells.Item(43, 8).Value = 18605.643   # H43: 19782.926 -> 18605.643
$ws.Cells.Item(43, 9).Value = 4850.533   # I43: 6138.5 -> 4850.533
$ws.Cells.Item(43, 11).Value = 4850.533   # K43: 6138.5 -> 4850.533
$ws.Cells.Item(43, 13).Value = -4699.533   # M43: -5987.5 -> -4699.533

$ws.Cells.Item(126, 8).Value = 7359.476   # H126: 7562.55 -> 7359.476
$ws.Cells.Item(126, 9).Value = 5456.375   # I126: 5764.7144 -> 5456.375
$ws.Cells.Item(126, 11).Value = 16369.125   # K126: 17294.1432 -> 16369.125
$ws.Cells.Item(126, 13).Value = -13899.125   # M126: -14824.1432 -> -13899.125

$ws.Cells.Item(133, 8).Value = 69714.28999999999   # H133: 68000 -> 69714.28999999999
$ws.Cells.Item(133, 10).Value = 69714.28999999999   # J133: 68000 -> 69714.28999999999
$ws.Cells.Item(133, 12).Value = 69714.28999999999   # L133: 68000 -> 69714.28999999999
$ws.Cells.Item(133, 14).Value = -79834.28999999999   # N133: -78120 -> -79834.28999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3602.342   # H46: 3737.25 -> 3602.342
$ws.Cells.Item(46, 9).Value = 1339.4   # I46: 1449.6666 -> 1339.4
$ws.Cells.Item(46, 11).Value = 1339.4   # K46: 1449.6666 -> 1339.4
$ws.Cells.Item(46, 13).Value = -1151.4   # M46: -1261.6666 -> -1151.4

$ws.Cells.Item(93, 8).Value = 5052473   # H93: 5052474 -> 5052473
$ws.Cells.Item(93, 9).Value = 1722.3793   # I93: 1723.6897 -> 1722.3793
$ws.Cells.Item(93, 11).Value = 1722.3793   # K93: 1723.6897 -> 1722.3793
$ws.Cells.Item(93, 13).Value = -474.3793000000001   # M93: -475.6896999999999 -> -474.3793000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 48653   # H96: 53348.3 -> 48653
$ws.Cells.Item(96, 9).Value = 84864.5   # I96: 101497.4 -> 84864.5
$ws.Cells.Item(96, 11).Value = 84864.5   # K96: 101497.4 -> 84864.5
$ws.Cells.Item(96, 13).Value = -83491.5   # M96: -100124.4 -> -83491.5

$ws.Cells.Item(122, 8).Value = 1456.5834   # H122: 1480.8182 -> 1456.5834
$ws.Cells.Item(122, 10).Value = 1538   # J122: 1625 -> 1538
$ws.Cells.Item(122, 12).Value = 4614   # L122: 4875 -> 4614
$ws.Cells.Item(122, 14).Value = -9514   # N122: -9775 -> -9514
